$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 5 (dropdown / select / None / id,name,class etc)
$ws.Range("E5").Value = "dropdown"
$ws.Range("F5").Value = "select"
$ws.Range("G5").Value = "None"
$ws.Range("H5").Value = "id,name,class etc"

# Fill row 6 (items in drop / option / value / id,name,class etc)
$ws.Range("E6").Value = "items in drop"
$ws.Range("F6").Value = "option"
$ws.Range("G6").Value = "value"
$ws.Range("H6").Value = "id,name,class etc"

# Update the active selection to F6
$ws.Range("F6").Select()
